# Weekly update: insert a new price record at the top of the
# "Femacal de La Calera - Alcachofa" time series (row 460), pushing the
# existing rows 460:483 down to 461:484.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 460 (shifts 460:483 -> 461:484).
$ws.Rows.Item(460).Insert()

# Populate the new row 460 with the latest weekly record. Columns that
# stay the same as the (now shifted-down) series are re-used; the
# observation-specific columns (date, quality, volume, prices) get the
# new values.
$ws.Range("A460").Value = 3
$ws.Range("B460").Value = "Femacal de La Calera"
$ws.Range("C460").Value = "Coquimbo"
$ws.Range("D460").Value = 45041
$ws.Range("E460").Value = 5
$ws.Range("F460").Value = 100112013
$ws.Range("G460").Value = "Alcachofa"
$ws.Range("H460").Value = "Argentina(o)"
$ws.Range("I460").Value = "Segunda"
$ws.Range("J460").Value = 75
$ws.Range("K460").Value = 12000
$ws.Range("L460").Value = 12000
$ws.Range("M460").Value = 12000
$ws.Range("N460").Value = "$/caja 50 unidades"
$ws.Range("O460").Value = "Provincia de Limarí"
$ws.Range("P460").Value = 240
$ws.Range("Q460").Value = 50
$ws.Range("R460").Value = "Hortaliza"
